$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly price-report row was inserted above the old row 579,
# shifting every subsequent row (old 579-656) down by one (new 580-657).
$ws.Rows(579).Insert()

# Populate the newly inserted row 579 with the new observation.
$ws.Cells.Item(579, 1).Value = 8
$ws.Cells.Item(579, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(579, 3).Value = "Coquimbo"
$ws.Cells.Item(579, 4).Value = 45127
$ws.Cells.Item(579, 5).Value = 4
$ws.Cells.Item(579, 6).Value = 100112043
$ws.Cells.Item(579, 7).Value = "Pepino dulce"
$ws.Cells.Item(579, 8).Value = "Sin especificar"
$ws.Cells.Item(579, 9).Value = "Primera"
$ws.Cells.Item(579, 10).Value = 400
$ws.Cells.Item(579, 11).Value = 14000
$ws.Cells.Item(579, 12).Value = 15000
$ws.Cells.Item(579, 13).Value = 14500
$ws.Cells.Item(579, 14).Value = "`$/bandeja 18 kilos"
$ws.Cells.Item(579, 15).Value = "Provincia de Limarí"
$ws.Cells.Item(579, 16).Value = 806
$ws.Cells.Item(579, 17).Value = 18
$ws.Cells.Item(579, 18).Value = "Hortaliza"
